# Consolidated user made libraries and added components for V3.
#
# 1. Rename the existing "Sheet1" to "PCB_Shield_v2".
# 2. Add a new sheet "PCB_Shield_v3" right after it, with a small parts
#    table for the two new V3 components (buck converter + XT30 right
#    angle connector).
# 3. Update each sheet's selection / active-tab state to match the saved
#    workbook (v2 selection -> C16, v3 becomes the active/selected sheet
#    with selection -> G6).

$wb = $excel.ActiveWorkbook

# --- Rename sheet1, add sheet2 right after it -----------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "PCB_Shield_v2"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "PCB_Shield_v3"

# --- Column widths for the new sheet --------------------------------------
# (nudged slightly from the literal "Excel" character widths in the source
# file so they land on the nearest width this engine's column grid can
# actually represent)
$ws2.Columns.Item(1).ColumnWidth = 11.0
$ws2.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws2.Columns.Item(3).ColumnWidth = 11.5
$ws2.Columns.Item(4).ColumnWidth = 15.5
$ws2.Columns.Item(5).ColumnWidth = 11.0
$ws2.Columns.Item(7).ColumnWidth = 16.5

# --- Page setup for the new sheet ------------------------------------------
$ws2.PageSetup.Orientation = 1

# --- Populate the new sheet's data -----------------------------------------
# Header row (write order matches original authoring order, leaving the
# "Need to Order (Y/N)" header for later / after row 2's "Y").
$ws2.Range("A1").Value = "Part #"
$ws2.Range("B1").Value = "Link"
$ws2.Range("C1").Value = "Spec"
$ws2.Range("D1").Value = "Description"
$ws2.Range("E1").Value = "Qty "
$ws2.Range("F1").Value = "Package"
$ws2.Range("H1").Value = "In Altium"
$ws2.Range("I1").Value = "Datasheet"

# Row 2 - TSR 1-2450 buck converter
$ws2.Range("A2").Value = "TSR 1-2450"
$ws2.Range("B2").Value = "https://www.digikey.com/en/products/detail/traco-power/TSR-1-2450/9383780"
$ws2.Range("C2").Value = "5V 1A Output"
$ws2.Range("D2").Value = "Buck Converter"
$ws2.Range("E2").Value = 1
$ws2.Range("F2").Value = "N/A"
$ws2.Range("G2").Value = "Y"
$ws2.Range("G1").Value = "Need to Order (Y/N)"
$ws2.Range("H2").Value = "N"
$ws2.Range("I2").Value = "https://www.tracopower.com/sites/default/files/products/datasheets/tsr1_datasheet.pdf"

# Row 3 - XT30 right angle connector (avionics battery connector)
$ws2.Range("A3").Value = "XT30 right angle connector"
$ws2.Range("B3").Value = "https://www.amazon.com/Amass-XT30PW-Upgrade-Connectors-Battery/dp/B099F2PXYN/ref=sr_1_1?crid=1HAEZSSILG0HP&keywords=XT30+right+angle+connector&qid=1643672653&sprefix=xt30+right+angle+connector%2Caps%2C116&sr=8-1"
$ws2.Range("C3").Value = "XT30"
$ws2.Range("D3").Value = "Avionics Battery Connector"
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Value = "N/A"
$ws2.Range("G3").Value = "N"
$ws2.Range("H3").Value = "N"
$ws2.Range("I3").Value = "https://www.amazon.com/Amass-XT30PW-Upgrade-Connectors-Battery/dp/B099F2PXYN/ref=sr_1_1?crid=1HAEZSSILG0HP&keywords=XT30+right+angle+connector&qid=1643672653&sprefix=xt30+right+angle+connector%2Caps%2C116&sr=8-1"

# --- Selections / active sheet ---------------------------------------------
# v2 keeps its own remembered selection (C16), not the active tab any more.
[void]$ws1.Range("C16").Select()

# v3 is the tab that was showing when the file was last saved.
$ws2.Activate()
[void]$ws2.Range("G6").Select()
